$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$oldText = [string]$ws1.Range("A1").Value()
$newText = $oldText -replace "1000 Bs = 3\.31 = 12868\.01 pesos", "1000 Bs = 3.33 = 12904.8 pesos"
$newText = $newText -replace "12868\.01 pesos = 3\.3 = 964\.11 Bs", "12904.8 pesos = 3.32 = 963.16 Bs"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 300
$ws2.Range("O10").Value = 3871.44
$ws2.Range("N12").Value = 3889
$ws2.Range("O12").Value = 290.26
